# Auto-generated edit script: apply Hyperion_Profits value updates per commit diff.
$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC (53 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 633.15
$ws.Range("J28").Value = 1657.5
$ws.Range("L28").Value = 1657.5
$ws.Range("N28").Value = -2627.5
$ws.Range("H69").Value = 52638036
$ws.Range("J69").Value = 66674116
$ws.Range("L69").Value = 200022348
$ws.Range("N69").Value = -200024096
$ws.Range("H72").Value = 52638036
$ws.Range("J72").Value = 66674116
$ws.Range("L72").Value = 600067044
$ws.Range("N72").Value = -600075780
$ws.Range("H100").Value = 3821
$ws.Range("I100").Value = 3875
$ws.Range("K100").Value = 3875
$ws.Range("M100").Value = -3334
$ws.Range("H129").Value = 44118680
$ws.Range("J129").Value = 1780.8572
$ws.Range("L129").Value = 5342.571599999999
$ws.Range("N129").Value = -15342.5716
$ws.Range("H132").Value = 3318.111
$ws.Range("I132").Value = 3270.5386
$ws.Range("K132").Value = 9811.6158
$ws.Range("M132").Value = -7281.6158
$ws.Range("H137").Value = 57988.688
$ws.Range("I137").Value = 70390.38
$ws.Range("J137").Value = 4248
$ws.Range("K137").Value = 211171.14
$ws.Range("L137").Value = 12744
$ws.Range("M137").Value = -208621.14
$ws.Range("N137").Value = -17844
$ws.Range("H138").Value = 3081.0896
$ws.Range("I138").Value = 1287.7858
$ws.Range("J138").Value = 3554.7925
$ws.Range("K138").Value = 3863.3574
$ws.Range("L138").Value = 10664.3775
$ws.Range("M138").Value = 1276.6426
$ws.Range("N138").Value = -20944.3775
$ws.Range("H139").Value = 88789
$ws.Range("I139").Value = 52472.668
$ws.Range("J139").Value = 110578.8
$ws.Range("K139").Value = 52472.668
$ws.Range("L139").Value = 110578.8
$ws.Range("M139").Value = -47332.668
$ws.Range("N139").Value = -120858.8
$ws.Range("H140").Value = 94115.8
$ws.Range("J140").Value = 94115.8
$ws.Range("L140").Value = 94115.8
$ws.Range("N140").Value = -104475.8
$ws.Range("H141").Value = 19530.924
$ws.Range("I141").Value = 15644.223
$ws.Range("K141").Value = 46932.669
$ws.Range("M141").Value = -41752.669

# ---- Worksheet: ARM (62 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9446.615
$ws.Range("I32").Value = 5623.709
$ws.Range("J32").Value = 18588.348
$ws.Range("K32").Value = 5623.709
$ws.Range("L32").Value = 18588.348
$ws.Range("M32").Value = -5336.709
$ws.Range("N32").Value = -19162.348
$ws.Range("H45").Value = 41804.54
$ws.Range("I45").Value = 64014.75
$ws.Range("J45").Value = 6268.2
$ws.Range("K45").Value = 64014.75
$ws.Range("L45").Value = 6268.2
$ws.Range("M45").Value = -63637.75
$ws.Range("N45").Value = -7022.2
$ws.Range("H61").Value = 3260.543
$ws.Range("I61").Value = 3104.9656
$ws.Range("J61").Value = 4012.5
$ws.Range("K61").Value = 3104.9656
$ws.Range("L61").Value = 4012.5
$ws.Range("M61").Value = -2892.9656
$ws.Range("N61").Value = -4436.5
$ws.Range("H63").Value = 4090.3
$ws.Range("I63").Value = 2192.8462
$ws.Range("K63").Value = 2192.8462
$ws.Range("M63").Value = -1506.8462
$ws.Range("H66").Value = 4090.3
$ws.Range("I66").Value = 2192.8462
$ws.Range("K66").Value = 10964.231
$ws.Range("M66").Value = -7532.231
$ws.Range("H88").Value = 2024.6666
$ws.Range("I88").Value = 1537
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 1537
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -1131
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 2024.6666
$ws.Range("I91").Value = 1537
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 1537
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -133
$ws.Range("N91").Value = -5808
$ws.Range("H104").Value = 10366.429
$ws.Range("J104").Value = 10366.429
$ws.Range("L104").Value = 10366.429
$ws.Range("N104").Value = -17354.429
$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 3260.543
$ws.Range("I136").Value = 3104.9656
$ws.Range("J136").Value = 4012.5
$ws.Range("K136").Value = 9314.8968
$ws.Range("L136").Value = 12037.5
$ws.Range("M136").Value = -6764.8968
$ws.Range("N136").Value = -17137.5

# ---- Worksheet: BSM (26 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 9747.333000000001
$ws.Range("J95").Value = 9747.333000000001
$ws.Range("L95").Value = 9747.333000000001
$ws.Range("N95").Value = -15239.333
$ws.Range("H107").Value = 2355.375
$ws.Range("I107").Value = 2169.611
$ws.Range("J107").Value = 2912.6667
$ws.Range("K107").Value = 2169.611
$ws.Range("L107").Value = 2912.6667
$ws.Range("M107").Value = -249.6109999999999
$ws.Range("N107").Value = -6752.6667
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H134").Value = 3563.3
$ws.Range("I134").Value = 1592
$ws.Range("J134").Value = 5534.6
$ws.Range("K134").Value = 4776
$ws.Range("L134").Value = 16603.8
$ws.Range("M134").Value = -2241
$ws.Range("N134").Value = -21673.8
$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200

# ---- Worksheet: CRP (47 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26158.719
$ws.Range("I31").Value = 1607.9565
$ws.Range("J31").Value = 61450.438
$ws.Range("K31").Value = 1607.9565
$ws.Range("L31").Value = 61450.438
$ws.Range("M31").Value = -1312.9565
$ws.Range("N31").Value = -62040.438
$ws.Range("H34").Value = 26158.719
$ws.Range("I34").Value = 1607.9565
$ws.Range("J34").Value = 61450.438
$ws.Range("K34").Value = 1607.9565
$ws.Range("L34").Value = 61450.438
$ws.Range("M34").Value = -1405.9565
$ws.Range("N34").Value = -61854.438
$ws.Range("H58").Value = 5260.2188
$ws.Range("I58").Value = 6599.1055
$ws.Range("J58").Value = 3303.3845
$ws.Range("K58").Value = 6599.1055
$ws.Range("L58").Value = 3303.3845
$ws.Range("M58").Value = -6396.1055
$ws.Range("N58").Value = -3709.3845
$ws.Range("H86").Value = 12569.917
$ws.Range("J86").Value = 13484.714
$ws.Range("L86").Value = 13484.714
$ws.Range("N86").Value = -15730.714
$ws.Range("H89").Value = 12569.917
$ws.Range("J89").Value = 13484.714
$ws.Range("L89").Value = 67423.57000000001
$ws.Range("N89").Value = -78655.57000000001
$ws.Range("H132").Value = 47679.586
$ws.Range("I132").Value = 26649.375
$ws.Range("J132").Value = 888888
$ws.Range("K132").Value = 79948.125
$ws.Range("L132").Value = 2666664
$ws.Range("M132").Value = -77418.125
$ws.Range("N132").Value = -2671724
$ws.Range("H136").Value = 5260.2188
$ws.Range("I136").Value = 6599.1055
$ws.Range("J136").Value = 3303.3845
$ws.Range("K136").Value = 19797.3165
$ws.Range("L136").Value = 9910.1535
$ws.Range("M136").Value = -17247.3165
$ws.Range("N136").Value = -15010.1535
$ws.Range("H141").Value = 163324.86
$ws.Range("J141").Value = 163324.86
$ws.Range("L141").Value = 163324.86
$ws.Range("N141").Value = -173684.86

# ---- Worksheet: CUL (21 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2784.1428
$ws.Range("I39").Value = 1899
$ws.Range("J39").Value = 2931.6667
$ws.Range("K39").Value = 5697
$ws.Range("L39").Value = 8795.000100000001
$ws.Range("M39").Value = -5403
$ws.Range("N39").Value = -9383.000100000001
$ws.Range("H129").Value = 1805.1666
$ws.Range("I129").Value = 1113.5
$ws.Range("J129").Value = 2496.8333
$ws.Range("K129").Value = 3340.5
$ws.Range("L129").Value = 7490.499899999999
$ws.Range("M129").Value = 1659.5
$ws.Range("N129").Value = -17490.4999
$ws.Range("H139").Value = 83335064
$ws.Range("I139").Value = 100001270
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 300003810
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -299998670
$ws.Range("N139").Value = -22280

# ---- Worksheet: GSM (26 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 48899
$ws.Range("J95").Value = 48899
$ws.Range("L95").Value = 48899
$ws.Range("N95").Value = -54391
$ws.Range("H97").Value = 995.5
$ws.Range("J97").Value = 1209.25
$ws.Range("L97").Value = 1209.25
$ws.Range("N97").Value = -2201.25
$ws.Range("H107").Value = 292.33334
$ws.Range("I107").Value = 81
$ws.Range("J107").Value = 334.6
$ws.Range("K107").Value = 81
$ws.Range("L107").Value = 334.6
$ws.Range("M107").Value = 1839
$ws.Range("N107").Value = -4174.6
$ws.Range("H127").Value = 48700.5
$ws.Range("I127").Value = 19324
$ws.Range("K127").Value = 19324
$ws.Range("M127").Value = -14364
$ws.Range("H132").Value = 2824.4666
$ws.Range("I132").Value = 2749.6206
$ws.Range("J132").Value = 4995
$ws.Range("K132").Value = 8248.861800000001
$ws.Range("L132").Value = 14985
$ws.Range("M132").Value = -5718.861800000001
$ws.Range("N132").Value = -20045

# ---- Worksheet: LTW (15 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7440.909
$ws.Range("I46").Value = 2989.75
$ws.Range("J46").Value = 9984.429
$ws.Range("K46").Value = 2989.75
$ws.Range("L46").Value = 9984.429
$ws.Range("M46").Value = -2801.75
$ws.Range("N46").Value = -10360.429
$ws.Range("H109").Value = 49995
$ws.Range("J109").Value = 49995
$ws.Range("L109").Value = 49995
$ws.Range("N109").Value = -52769
$ws.Range("H136").Value = 48167.71
$ws.Range("I136").Value = 68991.664
$ws.Range("K136").Value = 206974.992
$ws.Range("M136").Value = -204424.992

# ---- Worksheet: WVR (13 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 861718
$ws.Range("H122").Value = 2909.963
$ws.Range("J122").Value = 4789.1113
$ws.Range("L122").Value = 14367.3339
$ws.Range("N122").Value = -19267.3339
$ws.Range("H132").Value = 284029.22
$ws.Range("I132").Value = 5125.353
$ws.Range("K132").Value = 15376.059
$ws.Range("M132").Value = -12846.059
$ws.Range("H136").Value = 3478.5557
$ws.Range("I136").Value = 3215.0344
$ws.Range("K136").Value = 9645.1032
$ws.Range("M136").Value = -7095.1032

